$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 18.362432188124956
$ws.Range("C2").Value = 23.28494295750005
$ws.Range("D2").Value = 26.649947188124941
$ws.Range("E2").Value = 26.52725895750001

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 20.167593688124953
$ws.Range("C3").Value = 29.415951360000065
$ws.Range("D3").Value = 22.914452188124926
$ws.Range("E3").Value = 25.565617110000062

# Update the active selection to match the narrowed range
[void]$ws.Range("B1:E3").Select()
